$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data (row 2), pushing the existing
# data rows down by two.
$ws.Rows("2:3").Insert()
$ws.Range("A2:C3").Style = "Normal"
$ws.Range("A2:C3").ClearFormats()

# Populate the two newly inserted rows with the new sensor readings.
$ws.Range("A2").Value = -0.0397062413394451
$ws.Range("B2").Value = -0.0024434609804302
$ws.Range("C2").Value = 0.0332921557128429

$ws.Range("A3").Value = -0.0200058370828628
$ws.Range("B3").Value = -0.0035124751739203
$ws.Range("C3").Value = 0.0421497002243995

# The insert shifted the former last three data rows (20, 21, 22) down to
# rows 22, 23 and 24. Remove them so the sheet ends at row 21 again.
$ws.Rows("22:24").Delete()
